# Swap the data for rows 12 and 13 (two species records captured at the
# same locality get reordered), and move the "Aktivitet" (M) details from
# row 13 to row 12, per the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- capture current ("before") values for the columns that swap ---
$row12_A  = $ws.Range("A12").Value2
$row12_B  = $ws.Range("B12").Value2
$row12_D  = $ws.Range("D12").Value2
$row12_E  = $ws.Range("E12").Value2
$row12_F  = $ws.Range("F12").Value2
$row12_G  = $ws.Range("G12").Value2
$row12_H  = $ws.Range("H12").Value2
$row12_Q  = $ws.Range("Q12").Value2
$row12_R  = $ws.Range("R12").Value2
$row12_AX = $ws.Range("AX12").Value2

$row13_A  = $ws.Range("A13").Value2
$row13_B  = $ws.Range("B13").Value2
$row13_D  = $ws.Range("D13").Value2
$row13_E  = $ws.Range("E13").Value2
$row13_F  = $ws.Range("F13").Value2
$row13_G  = $ws.Range("G13").Value2
$row13_H  = $ws.Range("H13").Value2
$row13_Q  = $ws.Range("Q13").Value2
$row13_R  = $ws.Range("R13").Value2
$row13_AX = $ws.Range("AX13").Value2

$row13_M  = $ws.Range("M13").Value2

# --- row 12 gets row 13's identity/taxon data ---
$ws.Range("A12").Value = $row13_A
$ws.Range("B12").Value = $row13_B
$ws.Range("D12").Value = $row13_D
$ws.Range("E12").Value = $row13_E
$ws.Range("F12").Value = $row13_F
$ws.Range("G12").Value = $row13_G
$ws.Range("H12").Value = $row13_H
$ws.Range("Q12").Value = $row13_Q
$ws.Range("R12").Value = $row13_R
$ws.Range("AX12").Value = $row13_AX

# row 12 also gains the K/L/M/N "Aktivitet" block that used to live on row 13
$ws.Range("K12").Value = "'"
$ws.Range("K12").Style = "Normal"
$ws.Range("L12").Value = "'"
$ws.Range("L12").Style = "Normal"
$ws.Range("M12").Value = $row13_M
$ws.Range("N12").Value = "'"
$ws.Range("N12").Style = "Normal"

# --- row 13 gets row 12's (original) identity/taxon data ---
$ws.Range("A13").Value = $row12_A
$ws.Range("B13").Value = $row12_B
$ws.Range("D13").Value = $row12_D
$ws.Range("E13").Value = $row12_E
$ws.Range("F13").Value = $row12_F
$ws.Range("G13").Value = $row12_G
$ws.Range("H13").Value = $row12_H
$ws.Range("Q13").Value = $row12_Q
$ws.Range("R13").Value = $row12_R
$ws.Range("AX13").Value = $row12_AX

# row 13 loses the K/L/M/N "Aktivitet" block entirely (moved to row 12)
$ws.Range("K13").ClearContents()
$ws.Range("L13").ClearContents()
$ws.Range("M13").ClearContents()
$ws.Range("N13").ClearContents()
